# Update NATMI ligand-receptor TPM summary (Apln-Aplnr) values on Sheet1
# with the recomputed statistics (rows 2-13, columns G-T).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 20.419994
$ws.Range("H2").Value = 61.259982
$ws.Range("I2").Value = 0.9499654772891737
$ws.Range("J2").Value = 0.9499654772891736
$ws.Range("M2").Value = 231.368169
$ws.Range("N2").Value = 694.104507
$ws.Range("O2").Value = 0.9911105528536978
$ws.Range("P2").Value = 0.9911105528536978
$ws.Range("Q2").Value = 4724.536622770986
$ws.Range("R2").Value = 42520.82960493887
$ws.Range("S2").Value = 0.9415208093879998
$ws.Range("T2").Value = 0.9415208093879996
$ws.Range("G3").Value = 20.419994
$ws.Range("H3").Value = 61.259982
$ws.Range("I3").Value = 0.9499654772891737
$ws.Range("J3").Value = 0.9499654772891736
$ws.Range("O3").Value = 0.0032324887202399
$ws.Range("P3").Value = 0.0032324887202399
$ws.Range("Q3").Value = 15.408988732382
$ws.Range("R3").Value = 138.680898591438
$ws.Range("S3").Value = 0.003070752689954567
$ws.Range("T3").Value = 0.003070752689954567
$ws.Range("G4").Value = 20.419994
$ws.Range("H4").Value = 61.259982
$ws.Range("I4").Value = 0.9499654772891737
$ws.Range("J4").Value = 0.9499654772891736
$ws.Range("M4").Value = 1.150710666666667
$ws.Range("N4").Value = 3.452132
$ws.Range("O4").Value = 0.004929292953062385
$ws.Range("P4").Value = 0.004929292953062385
$ws.Range("Q4").Value = 23.49750490906933
$ws.Range("R4").Value = 211.477544181624
$ws.Range("S4").Value = 0.004682658132854069
$ws.Range("T4").Value = 0.004682658132854068
$ws.Range("G5").Value = 20.419994
$ws.Range("H5").Value = 61.259982
$ws.Range("I5").Value = 0.9499654772891737
$ws.Range("J5").Value = 0.9499654772891736
$ws.Range("M5").Value = 0.1698686666666667
$ws.Range("N5").Value = 0.509606
$ws.Range("O5").Value = 0.0007276654729999635
$ws.Range("P5").Value = 0.0007276654729999635
$ws.Range("Q5").Value = 3.468717154121333
$ws.Range("R5").Value = 31.218454387092
$ws.Range("S5").Value = 0.0006912570783652627
$ws.Range("T5").Value = 0.0006912570783652625
$ws.Range("I6").Value = 0.0191292957514309
$ws.Range("J6").Value = 0.0191292957514309
$ws.Range("M6").Value = 231.368169
$ws.Range("N6").Value = 694.104507
$ws.Range("O6").Value = 0.9911105528536978
$ws.Range("P6").Value = 0.9911105528536978
$ws.Range("Q6").Value = 95.137202883786
$ws.Range("R6").Value = 856.234825954074
$ws.Range("S6").Value = 0.01895924688790257
$ws.Range("T6").Value = 0.01895924688790257
$ws.Range("I7").Value = 0.0191292957514309
$ws.Range("J7").Value = 0.0191292957514309
$ws.Range("O7").Value = 0.0032324887202399
$ws.Range("P7").Value = 0.0032324887202399
$ws.Range("S7").Value = 0.00006183523274263343
$ws.Range("T7").Value = 0.00006183523274263343
$ws.Range("I8").Value = 0.0191292957514309
$ws.Range("J8").Value = 0.0191292957514309
$ws.Range("M8").Value = 1.150710666666667
$ws.Range("N8").Value = 3.452132
$ws.Range("O8").Value = 0.004929292953062385
$ws.Range("P8").Value = 0.004929292953062385
$ws.Range("Q8").Value = 0.4731653218693334
$ws.Range("R8").Value = 4.258487896824
$ws.Range("S8").Value = 0.00009429390274457457
$ws.Range("T8").Value = 0.00009429390274457454
$ws.Range("I9").Value = 0.0191292957514309
$ws.Range("J9").Value = 0.0191292957514309
$ws.Range("M9").Value = 0.1698686666666667
$ws.Range("N9").Value = 0.509606
$ws.Range("O9").Value = 0.0007276654729999635
$ws.Range("P9").Value = 0.0007276654729999635
$ws.Range("Q9").Value = 0.06984897652133333
$ws.Range("R9").Value = 0.628640788692
$ws.Range("S9").Value = 0.00001391972804112116
$ws.Range("T9").Value = 0.00001391972804112116
$ws.Range("G10").Value = 0.6643236666666666
$ws.Range("H10").Value = 1.992971
$ws.Range("I10").Value = 0.03090522695939548
$ws.Range("J10").Value = 0.03090522695939548
$ws.Range("M10").Value = 231.368169
$ws.Range("N10").Value = 694.104507
$ws.Range("O10").Value = 0.9911105528536978
$ws.Range("P10").Value = 0.9911105528536978
$ws.Range("Q10").Value = 153.703350380033
$ws.Range("R10").Value = 1383.330153420297
$ws.Range("S10").Value = 0.03063049657779546
$ws.Range("T10").Value = 0.03063049657779546
$ws.Range("G11").Value = 0.6643236666666666
$ws.Range("H11").Value = 1.992971
$ws.Range("I11").Value = 0.03090522695939548
$ws.Range("J11").Value = 0.03090522695939548
$ws.Range("O11").Value = 0.0032324887202399
$ws.Range("P11").Value = 0.0032324887202399
$ws.Range("Q11").Value = 0.5013006318376667
$ws.Range("R11").Value = 4.511705686539001
$ws.Range("S11").Value = 0.00009990079754269993
$ws.Range("T11").Value = 0.00009990079754269995
$ws.Range("G12").Value = 0.6643236666666666
$ws.Range("H12").Value = 1.992971
$ws.Range("I12").Value = 0.03090522695939548
$ws.Range("J12").Value = 0.03090522695939548
$ws.Range("M12").Value = 1.150710666666667
$ws.Range("N12").Value = 3.452132
$ws.Range("O12").Value = 0.004929292953062385
$ws.Range("P12").Value = 0.004929292953062385
$ws.Range("Q12").Value = 0.7644443293524444
$ws.Range("R12").Value = 6.879998964172
$ws.Range("S12").Value = 0.0001523409174637418
$ws.Range("T12").Value = 0.0001523409174637418
$ws.Range("G13").Value = 0.6643236666666666
$ws.Range("H13").Value = 1.992971
$ws.Range("I13").Value = 0.03090522695939548
$ws.Range("J13").Value = 0.03090522695939548
$ws.Range("M13").Value = 0.1698686666666667
$ws.Range("N13").Value = 0.509606
$ws.Range("O13").Value = 0.0007276654729999635
$ws.Range("P13").Value = 0.0007276654729999635
$ws.Range("Q13").Value = 0.1128477754917778
$ws.Range("R13").Value = 0.628640788692
$ws.Range("S13").Value = 0.00002248866659357973
$ws.Range("T13").Value = 0.00002248866659357973
